$d = $word.ActiveDocument

# The only semantic change in this revision is that the hidden "_GoBack"
# bookmark (Word's "last edit location" marker) moved earlier in the
# document: from right before "values indicate the significance..." (after
# "For B-E, p-") to right after "...across diets and sex" and before the
# following ". " -- the visible text itself is unchanged.

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Locate the new bookmark position: immediately after the phrase
# "across diets and sex" (i.e. right before the period that ends the
# sentence "A) Violin plot of calcium levels at 19 weeks across diets and sex.").
$findRange = $d.Content
$found = $findRange.Find.Execute("across diets and sex", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $insertPoint = $d.Range($findRange.End, $findRange.End)
    $d.Bookmarks.Add("_GoBack", $insertPoint)
}
